# Commit: "commit the android 8.1 for airplane mode code."
#
# 1) Windows sheet: rows 14-17, "Development Status" (col G) change from
#    "Discussing" to "Planning" (copy the look of the existing "Planning"
#    cells at G20/H20).
# 2) Linux sheet: rows 6-9, "Software / Operation System" (col F) change
#    from "Android 8.0" to "Android 8.1.0".
# 3) Linux sheet: rows 7 and 9, "Development Status" (col G) change from
#    "Planning" to "Processing" (copy the look of the existing "Processing"
#    cells at Windows!H10:H13).

$wb = $excel.ActiveWorkbook
$wsWindows = $wb.Worksheets.Item("Windows")
$wsLinux   = $wb.Worksheets.Item("Linux")

# --- 1) Windows!G14:G17 -> "Planning" ---
$wsWindows.Range("G20").Copy()
for ($r = 14; $r -le 17; $r++) {
    $cell = $wsWindows.Range("G$r")
    $cell.PasteSpecial(-4122)   # xlPasteFormats
    $cell.Value = "Planning"
}

# --- 2) Linux!F6:F9 -> "Android 8.1.0" ---
for ($r = 6; $r -le 9; $r++) {
    $wsLinux.Range("F$r").Value = "Android 8.1.0"
}

# --- 3) Linux!G7 and G9 -> "Processing" ---
$wsWindows.Range("H10").Copy()
foreach ($r in 7, 9) {
    $cell = $wsLinux.Range("G$r")
    $cell.PasteSpecial(-4122)   # xlPasteFormats
    $cell.Value = "Processing"
}
